$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column F ("Excel Tab" header + "Price & emissions parameters " values) is removed
# entirely. Deleting it shifts the old column G ("Source" header, with a couple of empty
# hyperlink-styled cells at rows 14-15) left into the F position - matching the target layout:
# Year | Category | Metric | Unit | Value | Source
$ws.Range("F:F").Delete() | Out-Null

$ws.Range("C5").Select() | Out-Null
